$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.0248213576677426
$ws.Range("B2").Value = 0.009190559660654225
$ws.Range("C2").Value = 0.07193214446306229
$ws.Range("D2").Value = 0.09174852557694076
$ws.Range("E2").Value = 0.09499251063071806
$ws.Range("F2").Value = 0.09334234154648872
$ws.Range("G2").Value = 0.0260220912641523

$ws.Range("A3").Value = 0.02554079966903386
$ws.Range("B3").Value = 0.004648160779900539
$ws.Range("C3").Value = -0.005592891294911867
$ws.Range("D3").Value = 0.05676937660284444
$ws.Range("E3").Value = 0.1065174908767818
$ws.Range("F3").Value = 0.07406513025462461
$ws.Range("G3").Value = 0.01197793127816921

$ws.Range("A4").Value = 0.06581173914934509
$ws.Range("B4").Value = 0.02931022851194365
$ws.Range("C4").Value = 0.08717288821935654
$ws.Range("D4").Value = 0.1316837502792451
$ws.Range("E4").Value = 0.1344966250960489
$ws.Range("F4").Value = 0.1330753250879459
$ws.Range("G4").Value = 0.02842338630346995

$ws.Range("A5").Value = -0.003432642124829049
$ws.Range("B5").Value = -0.0006012367432928545
$ws.Range("C5").Value = 0.001902705547299353
$ws.Range("D5").Value = 0.06803268640457041
$ws.Range("E5").Value = 0.07284204907249862
$ws.Range("F5").Value = 0.07035527363843443
$ws.Range("G5").Value = 0.02251122182635113
